$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric-looking values in columns G:K to be stored as text,
# matching the existing rows (t="str") rather than being auto-converted
# to numbers.
$ws.Range("A5:K7").NumberFormat = "@"

# Row 5
$ws.Cells.Item(5, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(5, 2).Value = " November 01 2020"
$ws.Cells.Item(5, 3).Value = "KKR won by 60 runs"
$ws.Cells.Item(5, 4).Value = "Rajasthan Royals"
$ws.Cells.Item(5, 5).Value = "Kolkata Knight Riders"
$ws.Cells.Item(5, 6).Value = "Kartik Tyagi "
$ws.Cells.Item(5, 7).Value = "2"
$ws.Cells.Item(5, 8).Value = "3"
$ws.Cells.Item(5, 9).Value = "0"
$ws.Cells.Item(5, 10).Value = "0"
$ws.Cells.Item(5, 11).Value = "66.66"

# Row 6
$ws.Cells.Item(6, 1).Value = " Abu Dhabi"
$ws.Cells.Item(6, 2).Value = " October 06 2020"
$ws.Cells.Item(6, 3).Value = "Mumbai won by 57 runs"
$ws.Cells.Item(6, 4).Value = "Rajasthan Royals"
$ws.Cells.Item(6, 5).Value = "Mumbai Indians"
$ws.Cells.Item(6, 6).Value = "Kartik Tyagi "
$ws.Cells.Item(6, 7).Value = "0"
$ws.Cells.Item(6, 8).Value = "0"
$ws.Cells.Item(6, 9).Value = "0"
$ws.Cells.Item(6, 10).Value = "0"
$ws.Cells.Item(6, 11).Value = "-"

# Row 7
$ws.Cells.Item(7, 1).Value = " Sharjah"
$ws.Cells.Item(7, 2).Value = " October 09 2020"
$ws.Cells.Item(7, 3).Value = "Capitals won by 46 runs"
$ws.Cells.Item(7, 4).Value = "Rajasthan Royals"
$ws.Cells.Item(7, 5).Value = "Delhi Capitals"
$ws.Cells.Item(7, 6).Value = "Kartik Tyagi "
$ws.Cells.Item(7, 7).Value = "2"
$ws.Cells.Item(7, 8).Value = "3"
$ws.Cells.Item(7, 9).Value = "0"
$ws.Cells.Item(7, 10).Value = "0"
$ws.Cells.Item(7, 11).Value = "66.66"
